# Apply updated "想去人数" (F) and "最低票价" (G) figures to the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1) ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 7287
$ws1.Range("G2").Value = 70

$ws1.Range("G3").Value = "不可售"

$ws1.Range("F8").Value = 138
$ws1.Range("F13").Value = 18
$ws1.Range("F16").Value = 1873
$ws1.Range("F17").Value = 51
$ws1.Range("F19").Value = 3814
$ws1.Range("F26").Value = 2497
$ws1.Range("F28").Value = 323
$ws1.Range("F31").Value = 45
$ws1.Range("F37").Value = 169
$ws1.Range("F38").Value = 34
$ws1.Range("F39").Value = 1490
$ws1.Range("F40").Value = 164

# --- Sheet "全部类型" (index 4) ---
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 7287
$ws4.Range("G2").Value = 70

$ws4.Range("G3").Value = "不可售"

$ws4.Range("F9").Value = 138
$ws4.Range("F14").Value = 18
$ws4.Range("F17").Value = 1873
$ws4.Range("F18").Value = 51
$ws4.Range("F20").Value = 3814
$ws4.Range("F27").Value = 2497
$ws4.Range("F29").Value = 323
$ws4.Range("F32").Value = 45
$ws4.Range("F38").Value = 169
$ws4.Range("F39").Value = 34
$ws4.Range("F40").Value = 1490
$ws4.Range("F41").Value = 164
